# Applies updated cached calculation values to the Weights.xlsx workbook
# (aerodynamic features testing in progress - recalculated weight estimation values)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 63840.73013384652
$ws.Range("C7").Value = 57456.65712046187
$ws.Range("C9").Value = 11811.504388009664
$ws.Range("C10").Value = 55669.22574583686
$ws.Range("C11").Value = 52029.22574583687
$ws.Range("C14").Value = 38029.22574583687
$ws.Range("C15").Value = 37250.91286142469
$ws.Range("C16").Value = 319.2255934121704
$ws.Range("C19").Value = 35474.33845483687
$ws.Range("C20").Value = 20763.048824670594

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 6384.511868243408
$ws.Range("C3").Value = 8158.333333333332
$ws.Range("D3").Value = 27.78319629904552
$ws.Range("C5").Value = 8158.333333333332
$ws.Range("D8").Value = 229.40654562188683
$ws.Range("C9").Value = 10245.0
$ws.Range("D9").Value = 60.46645712977179
$ws.Range("D10").Value = 0.17993751117817416
$ws.Range("D11").Value = 23.595979815620225
$ws.Range("D12").Value = 69.1906949649385
$ws.Range("C13").Value = 6465.0
$ws.Range("D13").Value = 1.2606779252293459
$ws.Range("C14").Value = 14396.0
$ws.Range("D14").Value = 125.48317392290824
$ws.Range("C15").Value = 7151.0
$ws.Range("D15").Value = 12.005430447535197

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 6767.582580338012
$ws.Range("C3").Value = 7183.75
$ws.Range("D3").Value = 6.149425067543128
$ws.Range("C5").Value = 7183.749999999999
$ws.Range("C8").Value = 6634.0
$ws.Range("D8").Value = -1.9738596278989216
$ws.Range("C9").Value = 6140.0
$ws.Range("D9").Value = -9.273364201884139
$ws.Range("C10").Value = 8397.0
$ws.Range("D10").Value = 24.076801432700144
$ws.Range("C11").Value = 7564.0
$ws.Range("D11").Value = 11.768122667255435

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 734.2188648479919
$ws.Range("D3").Value = 8.49966381140701
$ws.Range("D8").Value = 92.72182556804131
$ws.Range("D9").Value = -45.65653116491273
$ws.Range("D10").Value = -4.660580991074963
$ws.Range("D11").Value = 41.647136813260055
$ws.Range("D12").Value = -92.917643159337
$ws.Range("D13").Value = -30.94702080353572
$ws.Range("D14").Value = 0.37878829939678926
$ws.Range("D15").Value = 107.43133592941832

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 734.2188648479919
$ws.Range("D3").Value = -8.337958581419214
$ws.Range("D8").Value = 55.94804966459881
$ws.Range("D9").Value = -31.62801665359947
$ws.Range("D10").Value = 2.01317833954979
$ws.Range("D11").Value = -83.11130291841899
$ws.Range("D12").Value = -75.48414939770498
$ws.Range("D13").Value = 107.43133592941832
$ws.Range("D14").Value = -33.53480503377797

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 1213.0572549662475
$ws.Range("D3").Value = 14.586512244937348
$ws.Range("D10").Value = 16.07036636033946
$ws.Range("D11").Value = 13.267530809024443
$ws.Range("D12").Value = 14.421639565448274
$ws.Range("D17").Value = 16.07036636033946
$ws.Range("D18").Value = 13.267530809024443
$ws.Range("D19").Value = 14.421639565448274

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value = 5299.144850642028
$ws.Range("D3").Value = 21.7303328835241
$ws.Range("D11").Value = 30.47388201064785
$ws.Range("D12").Value = 11.489686855497178
$ws.Range("D13").Value = 23.227429784427315
$ws.Range("D18").Value = 30.47388201064785
$ws.Range("D19").Value = 11.489686855497178
$ws.Range("D20").Value = 23.227429784427315

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 2617.649865979797
$ws.Range("C3").Value = 2561.3404913372715
$ws.Range("D3").Value = -2.1511423423869114
$ws.Range("C5").Value = 2561.340491337271
$ws.Range("C9").Value = 2561.3404913372715
$ws.Range("D9").Value = -2.1511423423869114
$ws.Range("C11").Value = 395.373479134483
$ws.Range("C13").Value = 2165.967012202789

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 8682.936140811034
$ws.Range("C3").Value = 8260.622963499623
$ws.Range("D3").Value = -4.8637139610699185
$ws.Range("C4").Value = 8260.622963499623
$ws.Range("C8").Value = 8260.622963499623
$ws.Range("D8").Value = -4.863713961069908
$ws.Range("C21").Value = 1034.2688160330767
$ws.Range("C23").Value = 1034.2688160330765
$ws.Range("C26").Value = 532.1477925468873
$ws.Range("C28").Value = 532.1477925468872
$ws.Range("C36").Value = 785.4203668559742
$ws.Range("C38").Value = 785.4203668559741
$ws.Range("C41").Value = 3322.8792808498156
$ws.Range("C43").Value = 3322.879280849815
